$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("mult")
$rng = $ws2.Range("A1:T86")
$rng.AutoFilter(7, @("08"), 7)
$nm = $ws2.Names.Add("_FilterDatabase", "=mult!`$A`$1:`$T`$86", $false)
$nm.Visible = $false
